$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 9: "01-08-2021" must be kept as text (not auto-converted to a date),
# matching the other date-label strings already present in column A.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "01-08-2021"
$ws.Range("A9").Style = "Normal"

$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 3.5
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2.5
